# "Generate Report for Handoff"
# Update the handoff/handback status text + timestamps on each sheet, and
# narrow the "datetime" columns that now hold the shorter status text.

$wb = $excel.ActiveWorkbook

# Target "characters" width from the report generator is 17.2159881591797.
# The host's ColumnWidth setter quantizes to a pixel-ish grid (buckets of
# 1/6), so feed it a value whose bucket lands on the closest reachable
# width (17.166666666666668) instead of the unreachable exact figure.
$reportColWidth = 16.3333333333333

# --- Overview sheet ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "Ready for handoff"
$ws1.Range("F2").Value = "Ready for handoff"
$ws1.Range("G2").Value = "2016-08-13 07:15:16"
$ws1.Columns.Item(5).ColumnWidth = $reportColWidth
$ws1.Columns.Item(6).ColumnWidth = $reportColWidth

# --- zh-cn sheet ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("H2").Value = "2016-08-13 07:15:09"
$ws2.Columns.Item(3).ColumnWidth = $reportColWidth

# --- de-de sheet ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("H2").Value = "2016-08-13 07:15:16"
$ws3.Columns.Item(3).ColumnWidth = $reportColWidth
